$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.000005789925963522006
$ws.Range("E2").Value = 0.000005789925963522006

$ws.Range("D3").Value = 0.993155480577911
$ws.Range("E3").Value = 0.993155480577911

$ws.Range("D4").Value = 0.000000000000000000000000000000001163539726637278
$ws.Range("E4").Value = 0.000000000000000000000000000000001163539726637278

$ws.Range("D5").Value = 0.9750184524815226
$ws.Range("E5").Value = 0.9750184524815226

$ws.Range("D6").Value = 0.9999354764857952
$ws.Range("E6").Value = 0.9999354764857952

$ws.Range("D7").Value = 0.9999999999997229
$ws.Range("E7").Value = 0.0000000000002771116669464391

$ws.Range("D8").Value = 0.9864200542161468
$ws.Range("E8").Value = 0.01357994578385324

$ws.Range("D9").Value = 0.8503222436505193
$ws.Range("E9").Value = 0.1496777563494807

$ws.Range("D11").Value = 0.7347726701087769
$ws.Range("E11").Value = 0.2652273298912231
$ws.Range("F11").Value = 1.880642175674438

$wb.Save()
